$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2264705882352941
$ws.Range("C2").Value = 0.5176470588235295
$ws.Range("J2").Value = 0.01323529411764706
$ws.Range("P2").Value = 0.1426470588235294
$ws.Range("S2").Value = 0.1
$ws.Range("B3").Value = 0.01061007957559682
$ws.Range("C3").Value = 0.03978779840848806
$ws.Range("J3").Value = 0.03448275862068965
$ws.Range("P3").Value = 0.7214854111405835
$ws.Range("S3").Value = 0.1936339522546419
$ws.Range("J4").Value = 0.05617977528089887
$ws.Range("P4").Value = 0.6966292134831461
$ws.Range("S4").Value = 0.247191011235955
$ws.Range("B6").Value = 0.06318082788671024
$ws.Range("D6").Value = 0.0130718954248366
$ws.Range("E6").Value = 0.002178649237472767
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.2745098039215687
$ws.Range("O6").Value = 0.02396514161220044
$ws.Range("Q6").Value = 0.1503267973856209
$ws.Range("R6").Value = 0.06535947712418301
$ws.Range("S6").Value = 0.3485838779956427
$ws.Range("B7").Value = 0.1296703296703297
$ws.Range("D7").Value = 0.01318681318681319
$ws.Range("E7").Value = 0.002197802197802198
$ws.Range("F7").Value = 0.03516483516483516
$ws.Range("J7").Value = 0.1516483516483516
$ws.Range("O7").Value = 0.01538461538461539
$ws.Range("Q7").Value = 0.221978021978022
$ws.Range("R7").Value = 0.07252747252747253
$ws.Range("S7").Value = 0.3582417582417582
$ws.Range("B8").Value = 0.1083627797408716
$ws.Range("D8").Value = 0.01648998822143698
$ws.Range("E8").Value = 0.001177856301531213
$ws.Range("F8").Value = 0.05182567726737338
$ws.Range("J8").Value = 0.1036513545347468
$ws.Range("O8").Value = 0.0176678445229682
$ws.Range("Q8").Value = 0.2096584216725559
$ws.Range("R8").Value = 0.07302709069493522
$ws.Range("S8").Value = 0.4181389870435807
$ws.Range("B9").Value = 0.0783289817232376
$ws.Range("D9").Value = 0.01566579634464752
$ws.Range("E9").Value = 0.002610966057441253
$ws.Range("F9").Value = 0.05744125326370757
$ws.Range("J9").Value = 0.1227154046997389
$ws.Range("O9").Value = 0.02088772845953003
$ws.Range("Q9").Value = 0.2297650130548303
$ws.Range("R9").Value = 0.05744125326370757
$ws.Range("S9").Value = 0.4151436031331593
$ws.Range("B10").Value = 0.1090646442759119
$ws.Range("D10").Value = 0.02311303719754424
$ws.Range("E10").Value = 0.001083423618634886
$ws.Range("F10").Value = 0.06825568797399784
$ws.Range("J10").Value = 0.1180931744312026
$ws.Range("O10").Value = 0.02022390754785121
$ws.Range("Q10").Value = 0.2228241242325749
$ws.Range("R10").Value = 0.07836764174792343
$ws.Range("S10").Value = 0.358974358974359
$ws.Range("F11").Value = 0.001408450704225352
$ws.Range("G11").Value = 0.1366197183098592
$ws.Range("J11").Value = 0.1028169014084507
$ws.Range("K11").Value = 0.2042253521126761
$ws.Range("L11").Value = 0.5464788732394367
$ws.Range("S11").Value = 0.008450704225352112
$ws.Range("G12").Value = 0.7339901477832512
$ws.Range("J12").Value = 0.2142857142857143
$ws.Range("K12").Value = 0.009852216748768473
$ws.Range("L12").Value = 0.01970443349753695
$ws.Range("S12").Value = 0.02216748768472906
$ws.Range("F13").Value = 0.0101010101010101
$ws.Range("G13").Value = 0.7171717171717171
$ws.Range("J13").Value = 0.2323232323232323
$ws.Range("S13").Value = 0.04040404040404041
$ws.Range("F15").Value = 0.04158790170132325
$ws.Range("H15").Value = 0.1550094517958412
$ws.Range("I15").Value = 0.08128544423440454
$ws.Range("J15").Value = 0.3421550094517958
$ws.Range("K15").Value = 0.0888468809073724
$ws.Range("M15").Value = 0.01323251417769376
$ws.Range("O15").Value = 0.08695652173913043
$ws.Range("S15").Value = 0.1909262759924386
$ws.Range("F16").Value = 0.01428571428571429
$ws.Range("H16").Value = 0.1785714285714286
$ws.Range("I16").Value = 0.07857142857142857
$ws.Range("J16").Value = 0.4023809523809524
$ws.Range("K16").Value = 0.130952380952381
$ws.Range("M16").Value = 0.02857142857142857
$ws.Range("O16").Value = 0.06904761904761905
$ws.Range("S16").Value = 0.09761904761904762
$ws.Range("F17").Value = 0.02117420596727623
$ws.Range("H17").Value = 0.1453320500481232
$ws.Range("I17").Value = 0.08373435996150144
$ws.Range("J17").Value = 0.4456207892204042
$ws.Range("K17").Value = 0.1000962463907603
$ws.Range("M17").Value = 0.01539942252165544
$ws.Range("O17").Value = 0.0712223291626564
$ws.Range("S17").Value = 0.1174205967276227
$ws.Range("F18").Value = 0.0303030303030303
$ws.Range("H18").Value = 0.1322314049586777
$ws.Range("I18").Value = 0.08539944903581267
$ws.Range("J18").Value = 0.418732782369146
$ws.Range("K18").Value = 0.1184573002754821
$ws.Range("M18").Value = 0.01652892561983471
$ws.Range("N18").Value = 0.002754820936639119
$ws.Range("O18").Value = 0.08539944903581267
$ws.Range("S18").Value = 0.1101928374655647
$ws.Range("F19").Value = 0.01864801864801865
$ws.Range("H19").Value = 0.1961926961926962
$ws.Range("I19").Value = 0.07536907536907538
$ws.Range("J19").Value = 0.3776223776223776
$ws.Range("K19").Value = 0.1184926184926185
$ws.Range("M19").Value = 0.02408702408702409
$ws.Range("N19").Value = 0.000777000777000777
$ws.Range("O19").Value = 0.07459207459207459
$ws.Range("S19").Value = 0.1142191142191142
